$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Global)
$ws.Range("C4").Value = 57.55885974037317
$ws.Range("D4").Value = 19.25471565207577
$ws.Range("E4").Value = 9.751210713283157
$ws.Range("F4").Value = 2.197522793618164
$ws.Range("G4").Value = 88.76230889935026

# Row 5
$ws.Range("C5").Value = 66.68103388197095
$ws.Range("D5").Value = 25.63139893131246
$ws.Range("G5").Value = 92.31243281328341

# Row 6
$ws.Range("C6").Value = 87.2563188009497
$ws.Range("D6").Value = 3.926244939905627
$ws.Range("E6").Value = 0.1795562169405285
$ws.Range("G6").Value = 91.36211995779585

# Row 8
$ws.Range("C8").Value = 65.47666476161474
$ws.Range("G8").Value = 75.37525481742372

# Row 11
$ws.Range("C11").Value = 74.24552968637236
$ws.Range("E11").Value = 16.98929623880066
$ws.Range("G11").Value = 95.33920971283111

# Row 12
$ws.Range("C12").Value = 55.13793497726639
$ws.Range("G12").Value = 95.35073559458699

# Row 13
$ws.Range("C13").Value = 57.74197087134971
$ws.Range("G13").Value = 98.48406223524952

# Row 14
$ws.Range("D14").Value = 34.63020293295953
$ws.Range("G14").Value = 75.59669778319711

# Row 16
$ws.Range("C16").Value = 22.71058217014532
$ws.Range("G16").Value = 98.26540129188766

# Row 17
$ws.Range("C17").Value = 77.42533713214959
$ws.Range("G17").Value = 87.5552063088527

# Footnote text update (H21): add 'Deep Sea' to the list and remove ISSCAAP code 46
$ws.Range("H21").Value = "NOTE: Percent coverages are performed across FAO major fishing areas to be consistent with Fishstatj. `nThus, landings from areas such as 'Salmon', 'Tuna', 'Deep Sea', and 'Sharks' are added back into the FAO major fishing area from where they were reported. `nPercent coverage calculations do not include landings from ISSCAAP codes 61, 62, 63, 64, 71, 72, 73, 74, 81, 82, 83, 91, 92, 93, 94, `nexcept for stocks from these groups which are included in the assessment."
